# Update code for report co so
# Refresh the "last_edited_time" (column D) timestamps for the Notion export
# rows to reflect the latest sync pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-36: last_edited_time -> 2024-07-08T01:58:00.000Z
$ws.Range("D2:D36").Value = "2024-07-08T01:58:00.000Z"

# Rows 37-55: last_edited_time -> 2024-07-08T01:55:00.000Z
$ws.Range("D37:D55").Value = "2024-07-08T01:55:00.000Z"

# Rows 56-73: last_edited_time -> 2024-07-08T01:56:00.000Z
$ws.Range("D56:D73").Value = "2024-07-08T01:56:00.000Z"
